$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Overall)
$ws.Range("C2").Value = 5173.30597609562
$ws.Range("D2").Value = 4183.40169681997
$ws.Range("E2").Value = 6163.21025537126
$ws.Range("F2").Value = 49.9077596150444
$ws.Range("G2").Value = 39.8602772337274
$ws.Range("H2").Value = 60.6770473881395

# Row 3 (Acinobacter)
$ws.Range("B3").Value = 108
$ws.Range("C3").Value = 4944.17592592593
$ws.Range("D3").Value = 1238.24599630913
$ws.Range("E3").Value = 8650.10585554272
$ws.Range("F3").Value = 32.5396617829015
$ws.Range("G3").Value = 7.54678244745857
$ws.Range("H3").Value = 63.3406555338657

# Row 4 (klebsiella)
$ws.Range("C4").Value = 5352.8679245283
$ws.Range("D4").Value = 2224.87803180426
$ws.Range("E4").Value = 8480.85781725235
$ws.Range("F4").Value = 39.3094744946566
$ws.Range("G4").Value = 17.3977231146548
$ws.Range("H4").Value = 65.310954668377

# Row 5 (Clostridium) - unchanged

# Row 6 (Enterococcus)
$ws.Range("B6").Value = 480
$ws.Range("C6").Value = 4490.13125
$ws.Range("D6").Value = 2710.58975503441
$ws.Range("E6").Value = 6269.67274496559
$ws.Range("F6").Value = 37.0551130862287
$ws.Range("G6").Value = 22.7415476871315
$ws.Range("H6").Value = 53.0378618897625

# Row 7 (Escherichiacoli)
$ws.Range("C7").Value = 3091.62592592593
$ws.Range("D7").Value = 695.420024451337
$ws.Range("E7").Value = 5487.83182740052
$ws.Range("F7").Value = 27.8254842217933
$ws.Range("G7").Value = 11.1079252499772
$ws.Range("H7").Value = 47.0584063177734

# Row 8 (Pseudomonas)
$ws.Range("B8").Value = 167
$ws.Range("C8").Value = 5004.52694610778
$ws.Range("D8").Value = 1771.84582513497
$ws.Range("E8").Value = 8237.2080670806
$ws.Range("F8").Value = 29.5130207619676
$ws.Range("G8").Value = 7.9491671949834
$ws.Range("H8").Value = 55.3844553204609

# Row 9 (Candida)
$ws.Range("C9").Value = 4225.53103448276
$ws.Range("D9").Value = 2512.28954643103
$ws.Range("E9").Value = 5938.77252253448
$ws.Range("F9").Value = 43.477212421635
$ws.Range("G9").Value = 28.0575735987752
$ws.Range("H9").Value = 60.7535572146734

# Row 10 (Staphylococcus)
$ws.Range("C10").Value = 6404.95808383234
$ws.Range("D10").Value = 4547.11618993366
$ws.Range("E10").Value = 8262.79997773101
$ws.Range("F10").Value = 55.8864347601714
$ws.Range("G10").Value = 39.1337615202656
$ws.Range("H10").Value = 74.6562464545866

# Row 11 (Blood)
$ws.Range("C11").Value = 8924.96095444685
$ws.Range("D11").Value = 6825.86237875215
$ws.Range("E11").Value = 11024.0595301416
$ws.Range("F11").Value = 73.9079551834024
$ws.Range("G11").Value = 53.5882163254928
$ws.Range("H11").Value = 96.915997852189

# Row 12 (Urinary)
$ws.Range("C12").Value = 3849.92234548336
$ws.Range("D12").Value = 2597.48356986861
$ws.Range("E12").Value = 5102.36112109811
$ws.Range("F12").Value = 38.8878365922042
$ws.Range("G12").Value = 26.6189019043628
$ws.Range("H12").Value = 52.3455887165466

# Row 13 (Respiratory)
$ws.Range("C13").Value = 8659.11764705882
$ws.Range("D13").Value = 6201.0046028662
$ws.Range("E13").Value = 11117.2306912515
$ws.Range("F13").Value = 84.774609426636
$ws.Range("G13").Value = 58.4640329296849
$ws.Range("H13").Value = 115.453662623338

# Row 14 (Wound)
$ws.Range("C14").Value = 4280.10837438424
$ws.Range("D14").Value = 1228.00040083513
$ws.Range("E14").Value = 7332.21634793334
$ws.Range("F14").Value = 39.6894796126501
$ws.Range("G14").Value = 17.7051091679091
$ws.Range("H14").Value = 65.7799806006469
